# Insert a new data row at row 10, pushing the existing rows 10..121 down to 11..122
# (dimension grows from A1:R121 to A1:R122), then populate the new row 10 with the
# new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44515
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 100112009
$ws.Range("G10").Value = "Acelga"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 3000
$ws.Range("N10").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 750
$ws.Range("Q10").Value = 4
$ws.Range("R10").Value = "Hortaliza"
